# Auto-generated edit script implementing the target diff.
# Adds a new last data column (CR) "Signal_Value_123" to the
# Step1_Data and Step2_Sj sheets, updates the recomputed row 6
# ("signal segment 6") values, and refreshes the dependent
# Point_Exceeds_Cumulative_Value (F6) figures on the four
# Step3_DataPts_* sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Step1_Data (sheet1) and Step2_Sj (sheet2): add column CR with
# header "Signal_Value_123" and per-row values, and rewrite the
# BL..CQ values for row 6 (signal segment 6) to their recomputed
# figures.
# ---------------------------------------------------------------

$ws = $wb.Worksheets.Item("Step1_Data")

# New column CR (96): copy the header cell formatting (bold,
# border, centered) from the last existing header cell (CQ1)
# before writing the new header text.
$ws.Range("CQ1").Copy()
$ws.Range("CR1").PasteSpecial(-4122)

# New column CR (96) header + per-row values
$ws.Range("CR1").Value = "Signal_Value_123"
$ws.Range("CR2").Value = 0
$ws.Range("CR3").Value = 0
$ws.Range("CR4").Value = 0
$ws.Range("CR5").Value = 0
$ws.Range("CR6").Value = 0.01462714418853499
$ws.Range("CR7").Value = 0
$ws.Range("CR8").Value = 0
$ws.Range("CR9").Value = 0
$ws.Range("CR10").Value = 0
$ws.Range("CR11").Value = 0

# Row 6 (signal segment 6): updated BL..CQ values
$ws.Range("BL6").Value = 0
$ws.Range("BM6").Value = 0.2031639267808175
$ws.Range("BN6").Value = 0.0244918494262969
$ws.Range("BO6").Value = 0.1526365939864656
$ws.Range("BP6").Value = 0.05387594645449971
$ws.Range("BQ6").Value = 0.0007729427416179824
$ws.Range("BR6").Value = 0.002893938251909972
$ws.Range("BS6").Value = 0.009004887655079907
$ws.Range("BT6").Value = 0.002827230227144133
$ws.Range("BU6").Value = 0.02776297924570294
$ws.Range("BV6").Value = 0.02687875148740203
$ws.Range("BW6").Value = 0.001018492751184847
$ws.Range("BX6").Value = 0.009000257877667929
$ws.Range("BY6").Value = 0.1537940157154496
$ws.Range("BZ6").Value = 0.09148284286635103
$ws.Range("CA6").Value = 0.02083338071841738
$ws.Range("CB6").Value = 0.03329218548831696
$ws.Range("CC6").Value = 0.003636099670445284
$ws.Range("CD6").Value = 0.01628154143730474
$ws.Range("CE6").Value = 0.08440056620739658
$ws.Range("CF6").Value = 0.01420589417143166
$ws.Range("CG6").Value = 0.0003059567950599146
$ws.Range("CH6").Value = 0.003071281531790753
$ws.Range("CI6").Value = 0.001564138760766048
$ws.Range("CJ6").Value = 0.000008432939559665289
$ws.Range("CK6").Value = 0.01048468342774407
$ws.Range("CL6").Value = 0.01665309645894262
$ws.Range("CM6").Value = 0.002930800003927567
$ws.Range("CN6").Value = 0.0008269299190177809
$ws.Range("CO6").Value = 0.00848971949011087
$ws.Range("CP6").Value = 0.001575315632567778
$ws.Range("CQ6").Value = 0.007208177691075364

$ws = $wb.Worksheets.Item("Step2_Sj")

# New column CR (96): copy the header cell formatting (bold,
# border, centered) from the last existing header cell (CQ1)
# before writing the new header text.
$ws.Range("CQ1").Copy()
$ws.Range("CR1").PasteSpecial(-4122)

# New column CR (96) header + per-row values
$ws.Range("CR1").Value = "Signal_Value_123"
$ws.Range("CR2").Value = 1
$ws.Range("CR3").Value = 1
$ws.Range("CR4").Value = 1
$ws.Range("CR5").Value = 1
$ws.Range("CR6").Value = 1
$ws.Range("CR7").Value = 0.9999999999999999
$ws.Range("CR8").Value = 0.9999999999999997
$ws.Range("CR9").Value = 1
$ws.Range("CR10").Value = 1
$ws.Range("CR11").Value = 0.9999999999999998

# Row 6 (signal segment 6): updated BL..CQ values
$ws.Range("BL6").Value = 0
$ws.Range("BM6").Value = 0.2031639267808175
$ws.Range("BN6").Value = 0.2276557762071144
$ws.Range("BO6").Value = 0.3802923701935801
$ws.Range("BP6").Value = 0.4341683166480798
$ws.Range("BQ6").Value = 0.4349412593896978
$ws.Range("BR6").Value = 0.4378351976416078
$ws.Range("BS6").Value = 0.4468400852966877
$ws.Range("BT6").Value = 0.4496673155238318
$ws.Range("BU6").Value = 0.4774302947695347
$ws.Range("BV6").Value = 0.5043090462569367
$ws.Range("BW6").Value = 0.5053275390081216
$ws.Range("BX6").Value = 0.5143277968857894
$ws.Range("BY6").Value = 0.668121812601239
$ws.Range("BZ6").Value = 0.75960465546759
$ws.Range("CA6").Value = 0.7804380361860074
$ws.Range("CB6").Value = 0.8137302216743243
$ws.Range("CC6").Value = 0.8173663213447696
$ws.Range("CD6").Value = 0.8336478627820744
$ws.Range("CE6").Value = 0.9180484289894709
$ws.Range("CF6").Value = 0.9322543231609026
$ws.Range("CG6").Value = 0.9325602799559626
$ws.Range("CH6").Value = 0.9356315614877533
$ws.Range("CI6").Value = 0.9371957002485193
$ws.Range("CJ6").Value = 0.9372041331880789
$ws.Range("CK6").Value = 0.947688816615823
$ws.Range("CL6").Value = 0.9643419130747657
$ws.Range("CM6").Value = 0.9672727130786932
$ws.Range("CN6").Value = 0.968099642997711
$ws.Range("CO6").Value = 0.9765893624878219
$ws.Range("CP6").Value = 0.9781646781203897
$ws.Range("CQ6").Value = 0.9853728558114651

# ---------------------------------------------------------------
# Step3_DataPts_* sheets: refresh Point_Exceeds_Cumulative_Value
# (column F) for row 6, which is sourced from Step2_Sj row 6.
# ---------------------------------------------------------------

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F6").Value = 0.5043090462569367

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F6").Value = 0.75960465546759

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F6").Value = 0.8137302216743243

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F6").Value = 0.9180484289894709

